# "ultima aula desse projeto"
#
# 1) Bump the auto datetime fields cached on the slide master / all slide
#    layouts / the notes master from 13/02/2025 -> 14/02/2025 (one day
#    later - the date the deck was actually edited).
# 2) Slide 1 subtitle: "Nome do Aluno" -> "Nome do Aluno: Anthony Gabriel".
# 3) Slide 17 title textbox: merge the two runs "Intens" + "Entrada" into a
#    single run "IntensEntrada".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders (datetimeFigureOut / datetime1 fields).
# ---------------------------------------------------------------------

$master = $p.SlideMaster

# Slide master's own "Date Placeholder" shape.
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "2/13/2025") {
        $sh.TextFrame.TextRange.Text = "2/14/2025"
    }
}

# Every custom (slide) layout hanging off the master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "2/13/2025") {
            $sh.TextFrame.TextRange.Text = "2/14/2025"
        }
    }
}

# Notes master (pt-BR locale -> dd/mm/yyyy style datetimeFigureOut field).
$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $sh = $nm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "13/02/2025") {
        $sh.TextFrame.TextRange.Text = "14/02/2025"
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1 - subtitle "Nome do Aluno" -> "Nome do Aluno: Anthony Gabriel"
# ---------------------------------------------------------------------

$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "Nome do Aluno") {
        $tr = $sh.TextFrame.TextRange
        # Keep "Nome do " as-is and retype "Aluno" -> "Aluno: Anthony Gabriel",
        # the same way PowerPoint splits a run when you place the caret mid
        # run and type past the end of the original word.
        $sub = $tr.Characters(9, 5)
        $sub.Text = "Aluno: Anthony Gabriel"
    }
}

# ---------------------------------------------------------------------
# 3) Slide 17 - merge "Intens" + "Entrada" runs into "IntensEntrada"
# ---------------------------------------------------------------------

$s17 = $p.Slides.Item(17)
for ($i = 1; $i -le $s17.Shapes.Count; $i++) {
    $sh = $s17.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "Tabela IntensEntrada") {
        $tr = $sh.TextFrame.TextRange
        # Characters(8, 13) == "IntensEntrada" (chars 8..20), i.e. both the
        # "Intens" and "Entrada" runs without touching the "Tabela " prefix.
        $sub = $tr.Characters(8, 13)
        $sub.Text = "IntensEntrada"
    }
}
